$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = 2019
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = "shirts"
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 89.90000000000001
